$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Septiembre de 2020 a las 05:06"

# Row 30 (Pais index 34)
$ws.Range("B30").Value = 121604
$ws.Range("C30").Value = 835
$ws.Range("D30").Value = 73150
$ws.Range("E30").Value = 41400
$ws.Range("G30").Value = 46
$ws.Range("H30").Value = 7054

# Row 33 (Pais index 37)
$ws.Range("B33").Value = 106425
$ws.Range("C33").Value = 64
$ws.Range("D33").Value = 99893
$ws.Range("E33").Value = 4944

# Row 39 (Pais index 43)
$ws.Range("B39").Value = 88769
$ws.Range("C39").Value = 402
$ws.Range("D39").Value = 18576
$ws.Range("E39").Value = 60284
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = 9909

# Row 177 (Pais index 181)
$ws.Range("B177").Value = 466
$ws.Range("D177").Value = 374
$ws.Range("E177").Value = 91

# Row 198 (Pais index 202)
$ws.Range("B198").Value = 62
$ws.Range("C198").Value = 1
$ws.Range("E198").Value = 4
